$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4166.76
$ws.Range("I33").Value = 74.21429000000001
$ws.Range("J33").Value = 9375.454
$ws.Range("K33").Value = 74.21429000000001
$ws.Range("L33").Value = 9375.454
$ws.Range("M33").Value = 154.78571
$ws.Range("N33").Value = -9833.454

$ws.Range("H80").Value = 2418.4
$ws.Range("I80").Value = 807.4211
$ws.Range("K80").Value = 2422.2633
$ws.Range("M80").Value = -1424.2633

$ws.Range("H83").Value = 2418.4
$ws.Range("I83").Value = 807.4211
$ws.Range("K83").Value = 7266.7899
$ws.Range("M83").Value = -2274.7899

$ws.Range("H98").Value = 3274305.8
$ws.Range("I98").Value = 7033.567
$ws.Range("K98").Value = 7033.567
$ws.Range("M98").Value = -5535.567

$ws.Range("H122").Value = 3274305.8
$ws.Range("I122").Value = 7033.567
$ws.Range("K122").Value = 21100.701
$ws.Range("M122").Value = -18650.701

$ws.Range("H135").Value = 2914.1072
$ws.Range("I135").Value = 1503.8
$ws.Range("K135").Value = 13534.2
$ws.Range("M135").Value = -10999.2

$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200

$ws.Range("H139").Value = 70178.336
$ws.Range("J139").Value = 70178.336
$ws.Range("L139").Value = 70178.336
$ws.Range("N139").Value = -80458.336

$ws.Range("H141").Value = 3700.4167
$ws.Range("I141").Value = 3650.7896
$ws.Range("J141").Value = 3889
$ws.Range("K141").Value = 10952.3688
$ws.Range("L141").Value = 11667
$ws.Range("M141").Value = -5772.3688
$ws.Range("N141").Value = -22027

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3292.5652
$ws.Range("I61").Value = 4254.077
$ws.Range("K61").Value = 4254.077
$ws.Range("M61").Value = -4042.077

$ws.Range("H74").Value = 1021.4
$ws.Range("I74").Value = 649.8570999999999
$ws.Range("K74").Value = 649.8570999999999
$ws.Range("M74").Value = 224.1429000000001

$ws.Range("H77").Value = 1021.4
$ws.Range("I77").Value = 649.8570999999999
$ws.Range("K77").Value = 3249.2855
$ws.Range("M77").Value = 1118.7145

$ws.Range("H88").Value = 6500
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 7400
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 7400
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -8212

$ws.Range("H91").Value = 6500
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 7400
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 7400
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -10208

$ws.Range("H136").Value = 3292.5652
$ws.Range("I136").Value = 4254.077
$ws.Range("K136").Value = 12762.231
$ws.Range("M136").Value = -10212.231

$ws.Range("H139").Value = 76566.664
$ws.Range("J139").Value = 76566.664
$ws.Range("L139").Value = 76566.664
$ws.Range("N139").Value = -86846.664

$ws.Range("H141").Value = 59900
$ws.Range("J141").Value = 63058.332
$ws.Range("L141").Value = 63058.332
$ws.Range("N141").Value = -73418.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 367.66666
$ws.Range("J22").Value = 301.5
$ws.Range("L22").Value = 301.5
$ws.Range("N22").Value = -647.5

$ws.Range("H86").Value = 16184.429
$ws.Range("I86").Value = 24943.54
$ws.Range("K86").Value = 24943.54
$ws.Range("M86").Value = -23820.54

$ws.Range("H89").Value = 16184.429
$ws.Range("I89").Value = 24943.54
$ws.Range("K89").Value = 124717.7
$ws.Range("M89").Value = -119101.7

$ws.Range("H134").Value = 1413.6111
$ws.Range("I134").Value = 1071.7273
$ws.Range("J134").Value = 1950.8572
$ws.Range("K134").Value = 3215.1819
$ws.Range("L134").Value = 5852.571599999999
$ws.Range("M134").Value = -680.1819
$ws.Range("N134").Value = -10922.5716

$ws.Range("H138").Value = 67425
$ws.Range("J138").Value = 67425
$ws.Range("L138").Value = 67425
$ws.Range("N138").Value = -77705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 475
$ws.Range("I12").Value = 475
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 475
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -305
$ws.Range("N12").ClearContents()

$ws.Range("H50").Value = 8945.75
$ws.Range("J50").Value = 8945.75
$ws.Range("L50").Value = 8945.75
$ws.Range("N50").Value = -10195.75

$ws.Range("H58").Value = 1576.5454
$ws.Range("I58").Value = 760.4
$ws.Range("J58").Value = 2256.6667
$ws.Range("K58").Value = 760.4
$ws.Range("L58").Value = 2256.6667
$ws.Range("M58").Value = -557.4
$ws.Range("N58").Value = -2662.6667

$ws.Range("H59").Value = 11995
$ws.Range("J59").Value = 11995
$ws.Range("L59").Value = 11995
$ws.Range("N59").Value = -14285

$ws.Range("H74").Value = 17811.2
$ws.Range("J74").Value = 17811.2
$ws.Range("L74").Value = 17811.2
$ws.Range("N74").Value = -19559.2

$ws.Range("H77").Value = 17811.2
$ws.Range("J77").Value = 17811.2
$ws.Range("L77").Value = 53433.60000000001
$ws.Range("N77").Value = -62169.60000000001

$ws.Range("H105").Value = 1260
$ws.Range("I105").Value = 1440
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 1440
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = 307
$ws.Range("N105").Value = -4394

$ws.Range("H136").Value = 1576.5454
$ws.Range("I136").Value = 760.4
$ws.Range("J136").Value = 2256.6667
$ws.Range("K136").Value = 2281.2
$ws.Range("L136").Value = 6770.000100000001
$ws.Range("M136").Value = 268.8000000000002
$ws.Range("N136").Value = -11870.0001

$ws.Range("H138").Value = 46983.332
$ws.Range("J138").Value = 46983.332
$ws.Range("L138").Value = 46983.332
$ws.Range("N138").Value = -57263.332

$ws.Range("H140").Value = 89950
$ws.Range("J140").Value = 89950
$ws.Range("L140").Value = 89950
$ws.Range("N140").Value = -100310

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 689.2857
$ws.Range("J44").Value = 689.2857
$ws.Range("L44").Value = 2067.8571
$ws.Range("N44").Value = -2863.8571

$ws.Range("H121").Value = 64005.906
$ws.Range("I121").Value = 5165
$ws.Range("J121").Value = 77584.58
$ws.Range("K121").Value = 15495
$ws.Range("L121").Value = 232753.74
$ws.Range("M121").Value = -14185
$ws.Range("N121").Value = -235373.74

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 31990.4
$ws.Range("J136").Value = 31990.4
$ws.Range("L136").Value = 95971.20000000001
$ws.Range("N136").Value = -101071.2

$ws.Range("H138").Value = 68450
$ws.Range("J138").Value = 68450
$ws.Range("L138").Value = 68450
$ws.Range("N138").Value = -78730

$ws.Range("H140").Value = 89914.5
$ws.Range("J140").Value = 89914.5
$ws.Range("L140").Value = 89914.5
$ws.Range("N140").Value = -100274.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1244.2858
$ws.Range("I16").Value = 1233.6666
$ws.Range("J16").Value = 1252.25
$ws.Range("K16").Value = 1233.6666
$ws.Range("L16").Value = 1252.25
$ws.Range("M16").Value = -1063.6666
$ws.Range("N16").Value = -1592.25

$ws.Range("H22").Value = 623.5278
$ws.Range("I22").Value = 523.05554
$ws.Range("J22").Value = 724
$ws.Range("K22").Value = 523.05554
$ws.Range("L22").Value = 724
$ws.Range("M22").Value = -228.05554
$ws.Range("N22").Value = -1314

$ws.Range("H27").Value = 623.5278
$ws.Range("I27").Value = 523.05554
$ws.Range("J27").Value = 724
$ws.Range("K27").Value = 523.05554
$ws.Range("L27").Value = 724
$ws.Range("M27").Value = -416.05554
$ws.Range("N27").Value = -938

$ws.Range("H46").Value = 13682.5
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 15465.714
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 15465.714
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -15841.714

$ws.Range("H47").Value = 30926.545
$ws.Range("J47").Value = 30926.545
$ws.Range("L47").Value = 30926.545
$ws.Range("N47").Value = -31906.545

$ws.Range("H52").Value = 30926.545
$ws.Range("J52").Value = 30926.545
$ws.Range("L52").Value = 30926.545
$ws.Range("N52").Value = -31392.545

$ws.Range("H133").Value = 87406.664
$ws.Range("J133").Value = 87406.664
$ws.Range("L133").Value = 87406.664
$ws.Range("N133").Value = -92466.664

$ws.Range("H136").Value = 3080.68
$ws.Range("I136").Value = 4115.5
$ws.Range("J136").Value = 2883.5715
$ws.Range("K136").Value = 12346.5
$ws.Range("L136").Value = 8650.7145
$ws.Range("M136").Value = -9796.5
$ws.Range("N136").Value = -13750.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 64875
$ws.Range("J139").Value = 64875
$ws.Range("L139").Value = 64875
$ws.Range("N139").Value = -75155
